# Update crypto price/volume data per the Apr 22 2023 scrape refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $value) {
    $cell = $ws.Range($addr)
    # Force text interpretation so numeric-looking strings (and values
    # with embedded dots like "27.361.60") are not coerced into numbers
    # or have significant trailing zeros stripped, then restore the
    # original (default) cell style so no formatting side effects occur.
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextCell "D2" "27.361.60"
Set-TextCell "E2" "  -3.74%  "
Set-TextCell "D3" "1.857.80"
Set-TextCell "E3" "  -4.68%  "
Set-TextCell "E4" "  -1.06%  "
Set-TextCell "D5" "323.31"
Set-TextCell "E5" "  +0.49%  "
Set-TextCell "E6" "  -0.86%  "
Set-TextCell "D7" "0.4531"
Set-TextCell "E7" "  -5.66%  "
Set-TextCell "D8" "0.3864"
Set-TextCell "E8" "  -5.61%  "
Set-TextCell "D9" "48.47"
Set-TextCell "E9" "  -9.80%  "
Set-TextCell "D10" "0.07910"
Set-TextCell "E10" "  -7.13%  "
Set-TextCell "D11" "1.017"
Set-TextCell "E12" "  -4.49%  "
Set-TextCell "D13" "1.868.32"
Set-TextCell "E13" "  -5.25%  "
Set-TextCell "D14" "5.915"
Set-TextCell "E14" "  -3.97%  "
Set-TextCell "D15" "7.128"
Set-TextCell "E15" "  -5.88%  "
Set-TextCell "D16" "1.002"
Set-TextCell "E16" "  -1.07%  "
Set-TextCell "D17" "0.00001035"
Set-TextCell "E17" "  -3.57%  "
Set-TextCell "D18" "85.79"
Set-TextCell "E18" "  -5.03%  "
Set-TextCell "D19" "0.06513"
Set-TextCell "E19" "  -1.65%  "
Set-TextCell "E20" "  -7.44%  "
Set-TextCell "D21" "1.002"
Set-TextCell "E21" "  -0.90%  "
Set-TextCell "D23" "27.402.58"
Set-TextCell "E23" "  -3.76%  "
Set-TextCell "D25" "2.282"
Set-TextCell "E25" "  -0.63%  "
Set-TextCell "D26" "2.091.11"
Set-TextCell "E26" "  -5.14%  "
Set-TextCell "D27" "153.79"
Set-TextCell "E27" "  -1.63%  "
Set-TextCell "D28" "19.79"
Set-TextCell "E28" "  -2.53%  "
Set-TextCell "D29" "2.065"
Set-TextCell "E29" "  -5.10%  "
Set-TextCell "D30" "5.439"
Set-TextCell "E30" "  -6.89%  "
Set-TextCell "D31" "120.62"
Set-TextCell "D32" "1.484"
Set-TextCell "E32" "  +3.02%  "
Set-TextCell "D33" "0.09285"
Set-TextCell "E33" "  -4.16%  "
Set-TextCell "D34" "0.9334"
Set-TextCell "E34" "  -5.27%  "
Set-TextCell "D35" "3.616"
Set-TextCell "E35" "  -2.16%  "
Set-TextCell "E36" "  -6.60%  "
Set-TextCell "E37" "  -4.13%  "
Set-TextCell "B38" "Hedera"
Set-TextCell "C38" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextCell "D38" "0.05992"
Set-TextCell "E38" "  -3.22%  "
Set-TextCell "B39" "TrustWalletToken"
Set-TextCell "C39" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextCell "D39" "1.219"
Set-TextCell "E39" "  -2.27%  "
Set-TextCell "D40" "8.233"
Set-TextCell "E40" "  -10.08%  "
Set-TextCell "E41" "  -0.87%  "
Set-TextCell "D42" "0.5908"
Set-TextCell "E42" "  -5.14%  "
Set-TextCell "E43" "  -1.86%  "
Set-TextCell "D44" "10.10"
Set-TextCell "E44" "  -9.87%  "
Set-TextCell "D45" "1.274"
Set-TextCell "E45" "  -4.46%  "
Set-TextCell "D46" "0.5611"
Set-TextCell "E46" "  -5.91%  "
Set-TextCell "D47" "11.90"
Set-TextCell "E47" "  -7.88%  "
Set-TextCell "D48" "3.369"
Set-TextCell "E48" "  -1.14%  "
Set-TextCell "D49" "1.923"
Set-TextCell "E49" "  -6.83%  "
Set-TextCell "D50" "0.06776"
Set-TextCell "E50" "  -0.62%  "
Set-TextCell "D51" "108.18"
Set-TextCell "E51" "  -2.26%  "
